# Update the date line in the title paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-07 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-08 Wednesday", 2)

# Update the practice-table answers. Cells are addressed directly by
# (row, column) instead of text search/replace because one quotient
# string ("36÷5=7, 1") appears twice in the table but must become two
# different values depending on which cell it is in.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "79÷6=13, 1"
$t.Cell(1, 2).Range.Text  = "14÷6=2, 2"
$t.Cell(1, 3).Range.Text  = "57÷7=8, 1"
$t.Cell(1, 4).Range.Text  = "26÷3=8, 2"
$t.Cell(1, 5).Range.Text  = "76÷7=10, 6"

$t.Cell(5, 1).Range.Text  = "60÷8=7, 4"
$t.Cell(5, 2).Range.Text  = "27÷4=6, 3"
$t.Cell(5, 3).Range.Text  = "66÷3=22, 0"
$t.Cell(5, 4).Range.Text  = "20÷2=10, 0"
$t.Cell(5, 5).Range.Text  = "78÷3=26, 0"

$t.Cell(9, 1).Range.Text  = "89÷9=9, 8"
$t.Cell(9, 2).Range.Text  = "29÷4=7, 1"
$t.Cell(9, 3).Range.Text  = "18÷2=9, 0"
$t.Cell(9, 4).Range.Text  = "12÷3=4, 0"
$t.Cell(9, 5).Range.Text  = "12÷4=3, 0"

$t.Cell(13, 1).Range.Text = "19÷3=6, 1"
$t.Cell(13, 2).Range.Text = "64÷4=16, 0"
$t.Cell(13, 3).Range.Text = "97÷6=16, 1"
$t.Cell(13, 4).Range.Text = "54÷3=18, 0"
$t.Cell(13, 5).Range.Text = "58÷7=8, 2"

$t.Cell(17, 1).Range.Text = "33÷6=5, 3"
$t.Cell(17, 2).Range.Text = "66÷4=16, 2"
$t.Cell(17, 3).Range.Text = "10÷8=1, 2"
$t.Cell(17, 4).Range.Text = "74÷4=18, 2"
$t.Cell(17, 5).Range.Text = "73÷6=12, 1"
